$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the label text for the two "tax on millionaires" rows (A11, A12)
$ws.Range("A11").Value = "National tax on millionaires"
$ws.Range("A12").Value = "Global tax on millionaires"

# Overwrite the "United States" (column B) figures with the corrected values
$ws.Range("B2").Value = 0.409115351906478
$ws.Range("B3").Value = 0.530814200731826
$ws.Range("B4").Value = 0.45011341124838
$ws.Range("B5").Value = 0.308391082895318
$ws.Range("B6").Value = 0.344105383530252
$ws.Range("B7").Value = 0.390484602243659
$ws.Range("B8").Value = 0.415360140801995
$ws.Range("B9").Value = 0.44148547521565
$ws.Range("B10").Value = 0.335948883448662
$ws.Range("B11").Value = 0.619576204238963
$ws.Range("B12").Value = 0.58097404887994
